$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates as described by the source diff (crypto price/volume refresh).
# Values that look like plain numbers (e.g. "312.14") must be written back as
# TEXT (matching the inlineStr cells already in the sheet), so the number format
# is temporarily forced to Text and the style reset to Normal afterwards to avoid
# leaving the cell tagged with a different style than before.
$ws.Range("D2").Value = "45.367.56"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.368.56"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  -0.13%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.14"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "108.21"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -1.31%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -3.40%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.80"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -1.65%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.46"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  -4.18%  "
$ws.Range("D15").Value = "2.728.80"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "2.362.27"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "45.358.95"
$ws.Range("E18").Value = "  -0.39%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.34"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +10.29%  "
$ws.Range("E20").Value = "  -1.67%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -5.28%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "73.13"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.68%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "259.68"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.68%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.08%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.06"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.45%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.27"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("E29").Value = "  -1.69%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0972"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.60%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "22.28"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.20%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "36.81"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.89%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "166.63"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.87%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("E38").Value = "  +9.13%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("E41").Value = "  -3.41%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "98.38"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -7.43%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "69.99"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "12.80"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -7.84%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "1.814.95"
$ws.Range("E47").Value = "  +9.64%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.47%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "83.10"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.95%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "110.75"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -6.17%  "
$ws.Range("E51").Value = "  -0.59%  "
